$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.016449202858267
$ws.Range("D2").Value = 1.022531173563231
$ws.Range("E2").Value = 0.9926147277508489
$ws.Range("F2").Value = 1.014794091505186
$ws.Range("I2").Value = 1.026666337893759
$ws.Range("J2").Value = 1.021668949432451
$ws.Range("K2").Value = 1.025365280519094
$ws.Range("L2").Value = 0.9955398523336033
$ws.Range("M2").Value = 1.017651177085586
$ws.Range("N2").Value = 1.011270085431978
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.017842430221044
$ws.Range("D3").Value = 1.02358355699787
$ws.Range("E3").Value = 0.9936372048519304
$ws.Range("F3").Value = 1.01683739800815
$ws.Range("I3").Value = 1.026983458838593
$ws.Range("J3").Value = 1.02269529810395
$ws.Range("K3").Value = 1.026223678847742
$ws.Range("L3").Value = 0.9963617723202692
$ws.Range("M3").Value = 1.019496072775336
$ws.Range("N3").Value = 1.01161381095242
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.018741182803464
$ws.Range("D4").Value = 1.024261804396752
$ws.Range("E4").Value = 0.9942998659930995
$ws.Range("F4").Value = 1.018156249193296
$ws.Range("I4").Value = 1.027185667392901
$ws.Range("J4").Value = 1.023356301799418
$ws.Range("K4").Value = 1.026775820780073
$ws.Range("L4").Value = 0.9968940712668345
$ws.Range("M4").Value = 1.020686192550293
$ws.Range("N4").Value = 1.011835046414615
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.019118367943862
$ws.Range("D5").Value = 1.024546296128732
$ws.Range("E5").Value = 0.9945786998346017
$ws.Range("F5").Value = 1.018709923295942
$ws.Range("I5").Value = 1.02726996278961
$ws.Range("J5").Value = 1.023633450367238
$ws.Range("K5").Value = 1.027007157574632
$ws.Range("L5").Value = 0.997117960005301
$ws.Range("M5").Value = 1.021185663153862
$ws.Range("N5").Value = 1.01192777442281
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.01918166113436
$ws.Range("D6").Value = 1.024594025996272
$ws.Range("E6").Value = 0.9946255319796338
$ws.Range("F6").Value = 1.018802842989918
$ws.Range("I6").Value = 1.027284074620113
$ws.Range("J6").Value = 1.023679941814261
$ws.Range("K6").Value = 1.027045954262631
$ws.Range("L6").Value = 0.9971555583673453
$ws.Range("M6").Value = 1.021269476813427
$ws.Range("N6").Value = 1.011943327555966
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.018746225313369
$ws.Range("D7").Value = 1.024265608309705
$ws.Range("E7").Value = 0.9943035907982488
$ws.Range("F7").Value = 1.018163650418267
$ws.Range("I7").Value = 1.027186796550778
$ws.Range("J7").Value = 1.023360007960124
$ws.Range("K7").Value = 1.026778914985484
$ws.Range("L7").Value = 0.9968970624462087
$ws.Range("M7").Value = 1.020692869841063
$ws.Range("N7").Value = 1.011836286545376
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.016920627817269
$ws.Range("D8").Value = 1.022887397590317
$ws.Range("E8").Value = 0.9929600610674301
$ws.Range("F8").Value = 1.015485333011178
$ws.Range("I8").Value = 1.026774131125558
$ws.Range("J8").Value = 1.022016457938381
$ws.Range("K8").Value = 1.025656067103517
$ws.Range("L8").Value = 0.995817528259106
$ws.Range("M8").Value = 1.018275435294629
$ws.Range("N8").Value = 1.011386494565567
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.0136821119093
$ws.Range("D9").Value = 1.020437724434812
$ws.Range("E9").Value = 0.9906006454969559
$ws.Range("F9").Value = 1.01073957208988
$ws.Range("I9").Value = 1.02602393928668
$ws.Range("J9").Value = 1.019624761564318
$ws.Range("K9").Value = 1.023651918332633
$ws.Range("L9").Value = 0.9939188001724441
$ws.Range("M9").Value = 1.013986816149685
$ws.Range("N9").Value = 1.010584763527804
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.011507954971146
$ws.Range("D10").Value = 1.018790009931372
$ws.Range("E10").Value = 0.989033133672735
$ws.Range("F10").Value = 1.00755681684498
$ws.Range("I10").Value = 1.025508155898012
$ws.Range("J10").Value = 1.018013536937844
$ws.Range("K10").Value = 1.022298246128533
$ws.Range("L10").Value = 0.9926553831429383
$ws.Range("M10").Value = 1.011107215723856
$ws.Range("N10").Value = 1.010043963417883
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.01056278505421
$ws.Range("D11").Value = 1.01807297355176
$ws.Range("E11").Value = 0.988355674866747
$ws.Range("F11").Value = 1.006173874866675
$ws.Range("I11").Value = 1.025281062397113
$ws.Range("J11").Value = 1.017311771460966
$ws.Range("K11").Value = 1.021707833178488
$ws.Range("L11").Value = 0.9921088820399291
$ws.Range("M11").Value = 1.009855188781436
$ws.Range("N11").Value = 1.009808255930911
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.010211131619987
$ws.Range("D12").Value = 1.017806090337009
$ws.Range("E12").Value = 0.9881042295826724
$ws.Range("F12").Value = 1.005659444400562
$ws.Range("I12").Value = 1.02519614188865
$ws.Range("J12").Value = 1.017050479547869
$ws.Range("K12").Value = 1.021487879009612
$ws.Range("L12").Value = 0.9919059725120875
$ws.Range("M12").Value = 1.009389335243618
$ws.Range("N12").Value = 1.009720469383473
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.010286588708866
$ws.Range("D13").Value = 1.017863362444181
$ws.Range("E13").Value = 0.9881581567098651
$ws.Range("F13").Value = 1.005769825563486
$ws.Range("I13").Value = 1.02521438338165
$ws.Range("J13").Value = 1.01710655599023
$ws.Range("K13").Value = 1.021535089400099
$ws.Range("L13").Value = 0.9919494934313052
$ws.Range("M13").Value = 1.009489298749784
$ws.Range("N13").Value = 1.009739310551623
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.010533729079375
$ws.Range("D14").Value = 1.018050924047757
$ws.Range("E14").Value = 0.9883348863814464
$ws.Range("F14").Value = 1.006131367211374
$ws.Range("I14").Value = 1.025274054446787
$ws.Range("J14").Value = 1.017290185796311
$ws.Range("K14").Value = 1.021689664973846
$ws.Range("L14").Value = 0.9920921077337197
$ws.Range("M14").Value = 1.00981669753202
$ws.Range("N14").Value = 1.009801004266275
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.010685923815969
$ws.Range("D15").Value = 1.018166414607523
$ws.Range("E15").Value = 0.9884438009545853
$ws.Range("F15").Value = 1.00635402558339
$ws.Range("I15").Value = 1.025310744402493
$ws.Range("J15").Value = 1.017403243070246
$ws.Range("K15").Value = 1.021784817914008
$ws.Range("L15").Value = 0.9921799884222134
$ws.Range("M15").Value = 1.010018312702417
$ws.Range("N15").Value = 1.009838984654447
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.011570602806105
$ws.Range("D16").Value = 1.018837521450851
$ws.Range("E16").Value = 0.9890781214508737
$ws.Range("F16").Value = 1.007648495154472
$ws.Range("I16").Value = 1.025523147908403
$ws.Range("J16").Value = 1.01806002369952
$ws.Range("K16").Value = 1.022337339383393
$ws.Range("L16").Value = 0.9926916645766087
$ws.Range("M16").Value = 1.011190198391015
$ws.Range("N16").Value = 1.010059573873749
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.012124526978125
$ws.Range("D17").Value = 1.019257528697708
$ws.Range("E17").Value = 0.9894763578477731
$ws.Range("F17").Value = 1.008459182365839
$ws.Range("I17").Value = 1.025655374942263
$ws.Range("J17").Value = 1.018470901722414
$ws.Range("K17").Value = 1.022682774382972
$ws.Range("L17").Value = 0.9930127773692701
$ws.Range("M17").Value = 1.011923899070446
$ws.Range("N17").Value = 1.010197529719798
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.012447261032195
$ws.Range("D18").Value = 1.01950216840775
$ws.Range("E18").Value = 0.9897087662937551
$ws.Range("F18").Value = 1.008931582230535
$ws.Range("I18").Value = 1.025732138591392
$ws.Range("J18").Value = 1.018710165523831
$ws.Range("K18").Value = 1.022883849904875
$ws.Range("L18").Value = 0.9932001317071766
$ws.Range("M18").Value = 1.012351359778147
$ws.Range("N18").Value = 1.010277848912811
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.012557244225661
$ws.Range("D19").Value = 1.019585526228965
$ws.Range("E19").Value = 0.9897880325774039
$ws.Range("F19").Value = 1.009092581055447
$ws.Range("I19").Value = 1.025758251693838
$ws.Range("J19").Value = 1.018791681694611
$ws.Range("K19").Value = 1.022952342008016
$ws.Range("L19").Value = 0.993264023964098
$ws.Range("M19").Value = 1.012497029639428
$ws.Range("N19").Value = 1.010305210657848
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.012065133539757
$ws.Range("D20").Value = 1.019212501475963
$ws.Range("E20").Value = 0.9894336180355766
$ws.Range("F20").Value = 1.008372251087716
$ws.Range("I20").Value = 1.025641225716924
$ws.Range("J20").Value = 1.018426859264225
$ws.Range("K20").Value = 1.022645755031794
$ws.Range("L20").Value = 0.9929783193490043
$ws.Range("M20").Value = 1.011845231241175
$ws.Range("N20").Value = 1.010182743707503
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.010460968380866
$ws.Range("D21").Value = 1.017995706914544
$ws.Range("E21").Value = 0.9882828385668255
$ws.Range("F21").Value = 1.006024922953528
$ws.Range("I21").Value = 1.025256498510532
$ws.Range("J21").Value = 1.017236128726896
$ws.Range("K21").Value = 1.0216441642862
$ws.Range("L21").Value = 0.9920501090198107
$ws.Range("M21").Value = 1.009720308879404
$ws.Range("N21").Value = 1.009782843499347
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.009449031580144
$ws.Range("D22").Value = 1.01722750836108
$ws.Range("E22").Value = 0.9875604150241496
$ws.Range("F22").Value = 1.004544746076026
$ws.Range("I22").Value = 1.025011317895678
$ws.Range("J22").Value = 1.016483848343658
$ws.Range("K22").Value = 1.021010668344455
$ws.Range("L22").Value = 0.991467000034148
$ws.Range("M22").Value = 1.008379674897886
$ws.Range("N22").Value = 1.009530053156474
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.0099857987146
$ws.Range("D23").Value = 1.017635046499837
$ws.Range("E23").Value = 0.9879432794636459
$ws.Range("F23").Value = 1.005329833588709
$ws.Range("I23").Value = 1.025141605596475
$ws.Range("J23").Value = 1.01688299293072
$ws.Range("K23").Value = 1.021346855231878
$ws.Range("L23").Value = 0.9917760702887607
$ws.Range("M23").Value = 1.009090815029372
$ws.Range("N23").Value = 1.009664191862722
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.012091971992022
$ws.Range("D24").Value = 1.019232848400503
$ws.Range("E24").Value = 0.9894529299347241
$ws.Range("F24").Value = 1.008411533012014
$ws.Range("I24").Value = 1.025647620261847
$ws.Range("J24").Value = 1.01844676137394
$ws.Range("K24").Value = 1.022662483752709
$ws.Range("L24").Value = 0.9929938892766438
$ws.Range("M24").Value = 1.011880779374433
$ws.Range("N24").Value = 1.010189425327675
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.01452196831937
$ws.Range("D25").Value = 1.021073565690156
$ws.Range("E25").Value = 0.9912096547607046
$ws.Range("F25").Value = 1.011969703825315
$ws.Range("I25").Value = 1.026220628111268
$ws.Range("J25").Value = 1.020245989694683
$ws.Range("K25").Value = 1.024173106746793
$ws.Range("L25").Value = 0.9944092447426411
$ws.Range("M25").Value = 1.015099057761346
$ws.Range("N25").Value = 1.010793130469202
